# Add two new columns, "I0" (I) and "IF" (J), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold/bordered/centered) from the existing
# header cell H1 onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data rows (plain/default formatting, matching columns A-H).
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 7

$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 7

$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 6

$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 6

$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 4

$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 5
